$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51.
# NumberFormat "@" + ClearFormats() pattern forces text-typed values
# without altering cell style indices (values like "4.70", "233.28",
# "0.0978" would otherwise be auto-coerced to numeric by Excel).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "37.784.95"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +0.03%  "
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.079.77"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -0.78%  "
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  -0.14%  "
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "233.28"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +0.19%  "
$c.ClearFormats()

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +0.36%  "
$c.ClearFormats()

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "58.79"
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +1.92%  "
$c.ClearFormats()

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -0.17%  "
$c.ClearFormats()

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +1.14%  "
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0786"
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +1.12%  "
$c.ClearFormats()

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +2.41%  "
$c.ClearFormats()

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "14.85"
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +2.21%  "
$c.ClearFormats()

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "21.06"
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -1.53%  "
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.775"
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -0.27%  "
$c.ClearFormats()

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +2.39%  "
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.041.17"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -2.84%  "
$c.ClearFormats()

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "37.709.57"
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -0.17%  "
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.14"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -0.64%  "
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "71.53"
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +1.32%  "
$c.ClearFormats()

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +1.60%  "
$c.ClearFormats()

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "228.24"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +0.45%  "
$c.ClearFormats()

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "
$c.ClearFormats()

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +0.99%  "
$c.ClearFormats()

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.34"
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -2.76%  "
$c.ClearFormats()

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "170.58"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +1.60%  "
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.23"
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +2.99%  "
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -1.79%  "
$c.ClearFormats()

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +0.15%  "
$c.ClearFormats()

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +0.37%  "
$c.ClearFormats()

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +2.06%  "
$c.ClearFormats()

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.70"
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +1.19%  "
$c.ClearFormats()

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0634"
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +1.64%  "
$c.ClearFormats()

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.70"
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +2.23%  "
$c.ClearFormats()

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.49"
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -3.85%  "
$c.ClearFormats()

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -0.32%  "
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.40"
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  -1.33%  "
$c.ClearFormats()

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -0.11%  "
$c.ClearFormats()

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -1.32%  "
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0978"
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -0.44%  "
$c.ClearFormats()

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "99.60"
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +2.88%  "
$c.ClearFormats()

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -2.42%  "
$c.ClearFormats()

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0216"
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +0.54%  "
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "16.73"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +5.54%  "
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.440.29"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -0.85%  "
$c.ClearFormats()

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -0.47%  "
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.20"
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +2.25%  "
$c.ClearFormats()

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +0.42%  "
$c.ClearFormats()

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.42"
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +1.93%  "
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.96"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -1.57%  "
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.269.61"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -1.10%  "
$c.ClearFormats()

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "46.49"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +0.71%  "
$c.ClearFormats()

